$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# RoomStock data refresh: update cell values that differ in the new export
# (room 2 bumped to floor 2/price 2800; rows 45-54 renumbered/retyped; the
# last block of sample rows (8001-8004, 1200-9999) was regenerated).

# Row 2
$ws.Range("F2").Value = 2
$ws.Range("H2").Value = 2800

# Row 45
$ws.Range("A45").Value = 2052
$ws.Range("B45").Value = 7004

# Row 46
$ws.Range("A46").Value = 2054
$ws.Range("B46").Value = 7106

# Row 47
$ws.Range("A47").Value = 2057
$ws.Range("B47").Value = 8001
$ws.Range("C47").Value = "Delux"
$ws.Range("D47").Value = "Premium"
$ws.Range("F47").Value = 8
$ws.Range("H47").Value = 6500

# Row 48
$ws.Range("A48").Value = 2059
$ws.Range("B48").Value = 8003
$ws.Range("D48").Value = "Premium"
$ws.Range("F48").Value = 8
$ws.Range("H48").Value = 7500

# Row 49
$ws.Range("A49").Value = 2060
$ws.Range("B49").Value = 8004
$ws.Range("C49").Value = "Junior Suite"
$ws.Range("H49").Value = 7500

# Row 50
$ws.Range("A50").Value = 2052
$ws.Range("B50").Value = 1200
$ws.Range("C50").Value = "Superior"
$ws.Range("D50").Value = "Normal"
$ws.Range("F50").Value = 1
$ws.Range("H50").Value = 1200
$ws.Range("J50").Value = "24-04-2020 03:28:54"

# Row 51
$ws.Range("A51").Value = 2053
$ws.Range("B51").Value = 1200
$ws.Range("C51").Value = "Superior"
$ws.Range("D51").Value = "Normal"
$ws.Range("F51").Value = 1
$ws.Range("H51").Value = 1200
$ws.Range("J51").Value = "24-04-2020 03:30:52"

# Row 52
$ws.Range("A52").Value = 2054
$ws.Range("B52").Value = 1300
$ws.Range("C52").Value = "Superior"
$ws.Range("D52").Value = "Normal"
$ws.Range("F52").Value = 1
$ws.Range("H52").Value = 1300
$ws.Range("J52").Value = "24-04-2020 03:31:01"

# Row 53
$ws.Range("A53").Value = 2055
$ws.Range("B53").Value = 1102
$ws.Range("H53").Value = 1102
$ws.Range("J53").Value = "24-04-2020 03:33:59"

# Row 54
$ws.Range("A54").Value = 2056
$ws.Range("B54").Value = 9999
$ws.Range("H54").Value = 9999
$ws.Range("J54").Value = "24-04-2020 03:34:24"

# The refreshed export has 4 fewer sample rows; drop the now-unused tail rows
$ws.Range("A55:J58").EntireRow.Delete()
